# Add the IUSSI 2024 congress entry to the "conferences" sheet.
#
# A new row is inserted right under the header row (row 2), pushing the
# existing entries down by one, and the four data columns are populated
# with the new conference's title, date, event name and location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conferences")

# Insert a fresh row at position 2 (below the header), shifting the
# existing conference rows down.
$ws.Rows.Item(2).Insert()

# Populate the new entry: Title / Year(date) / Unit(event) / Where.
$ws.Range("A2").Value = "Honey bee (\textit{Apis mellifera}) CHC variations: New insight from desaturases and elongases`nexpression assays"
$ws.Range("B2").Value = "July 8 2024"
$ws.Range("C2").Value = "European meeting of the International Union for the Study of Social Insects (IUSSI)"
$ws.Range("D2").Value = "Lausanne, Switzerland"

# The row-insert carries the row's style into column E too (like the
# header row above it); the other data rows never populate column E, so
# drop that empty cell to match the rest of the table.
$ws.Range("E2").Clear()

# Match the row height used by the other wrapped, multi-line entries.
$ws.Rows.Item(2).RowHeight = 60

# Bring the conferences sheet to the foreground (it becomes the active
# tab after adding the new congress) and focus/zoom on the new entry.
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("E2").Select()
